$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows scraped for the newsletter feed (rows 337-350).
$newRows = @(
    @{ Row = 337; A = 'We are deeply saddened by the passing of Brent Schroeder, a friend, mentor, former colleague and a leader within the air-conditi'; B = 'We are deeply saddened by the passing of Brent Schroeder, a friend, mentor, former colleague and a leader within the air-conditioning and refrigeration industry. Brent spent 35 years at Copeland...'; C = 'https://www.linkedin.com/feed/update/urn:li:activity:7335787631324057601'; D = '2025-06-03' }
    @{ Row = 338; A = 'The next innovation in demanding data center cooling, health care and large chiller applications is here! The Copeland oil-free '; B = 'The next innovation in demanding data center cooling, health care and large chiller applications is here! The Copeland oil-free centrifugal compressor with frictionless Aero-lift bearing technology...'; C = 'https://www.linkedin.com/feed/update/urn:li:activity:7335742157892255744'; D = '2025-06-03' }
    @{ Row = 339; A = 'Drive performance, efficiency, and innovation with Danfoss at iVT Expo Cologne 2025'; B = ' June 3, 2025 June 2, 2025 — NORDBORG, DENMARK — Danfoss Power Solutions will highlight its latest technologies for off- and on-highway machinery at iVT Expo 2025, June 11-12, in Cologne, Germany.'; C = 'https://www.ejarn.com/article/detail/89014'; D = '2025-06-03' }
    @{ Row = 340; A = 'DOE Coolerchips’ Peter De Bock joins Eaton as VP of data center energy & cooling'; B = ' As ARPA-E faces dramatic cuts to its budget The head of the Department of Energy’s data center cooling research effort has left for equipment supplier Eaton.'; C = 'https://www.datacenterdynamics.com/en/news/doe-coolerchips-peter-de-bock-joins-eaton-as-vp-of-data-center-energy-cooling/'; D = '2025-06-03' }
    @{ Row = 341; A = 'TX2-W-G04-Y and TX2-W-G04-Y/H are the Climaveneta branded water source chillers and reversible heat pumps, dedicated to process '; B = 'TX2-W-G04-Y and TX2-W-G04-Y/H are the Climaveneta branded water source chillers and reversible heat pumps, dedicated to process applications. The range, with capacity from 191 to 2069 kW, has...'; C = 'https://www.linkedin.com/feed/update/urn:li:activity:7335636602922110976'; D = '2025-06-03' }
    @{ Row = 342; A = 'Romanian hazelnut farm secures €4.6M for expansion'; B = ' Dorin Bob, an entrepreneur from Transylvania, Romania, has secured €4.6 million in European funding to expand and modernize a major hazelnut plantation.'; C = 'https://www.freshplaza.com/europe/article/9737382/romanian-hazelnut-farm-secures-eur4-6m-for-expansion/'; D = '2025-06-03' }
    @{ Row = 343; A = 'On the occasion of China Refrigeration Expo, we interviewed Alvise Dina, Sales Director of #FrascoldChina, right at our booth. I'; B = 'On the occasion of China Refrigeration Expo, we interviewed Alvise Dina, Sales Director of #FrascoldChina, right at our booth.In this short video, Alvise shares insights into #Frascold’s presence in...'; C = 'https://www.linkedin.com/feed/update/urn:li:activity:7335597810764894210'; D = '2025-06-03' }
    @{ Row = 344; A = 'At Danfoss, we recognize waste heat as a pivotal element in our journey toward decarbonization. Discover how recovering and reus'; B = 'At Danfoss, we recognize waste heat as a pivotal element in our journey toward decarbonization. Discover how recovering and reusing excess heat can significantly enhance energy efficiency and drive...'; C = 'https://www.linkedin.com/feed/update/urn:li:activity:7335587321536294913'; D = '2025-06-03' }
    @{ Row = 345; A = 'Food retail leaders globally are adopting solutions to meet sustainability and zero-emissions goals, supported by low-GWP refrig'; B = 'Food retail leaders globally are adopting solutions to meet sustainability and zero-emissions goals, supported by low-GWP refrigerants that reduce energy waste. Copeland’s CO₂ scroll refrigeration...'; C = 'https://www.linkedin.com/feed/update/urn:li:activity:7335545929338392577'; D = '2025-06-03' }
    @{ Row = 346; A = 'Explore Scout, the new AI feature in Copeland Mobile. Scout provides tailored results for Copeland products and continuously imp'; B = 'Explore Scout, the new AI feature in Copeland Mobile. Scout provides tailored results for Copeland products and continuously improves through advanced learning. Available 24/7 on mobile and desktop,...'; C = 'https://www.linkedin.com/feed/update/urn:li:activity:7335341999199264768'; D = '2025-06-02' }
    @{ Row = 347; A = 'Danfoss Opens Danfoss Nanjing Park'; B = ' On April 23, the opening ceremony for Danfoss Nanjing Park was held in the Nanjing Economic Development Zone.'; C = 'https://www.ejarn.com/article/detail/88976'; D = '2025-06-01' }
    @{ Row = 348; A = 'Hanbell Establishes a New Company for Reciprocating Compressors'; B = ' Shanghai Hanbell Reciprocating Machinery was officially established with a registered capital of RMB 10 million (about US$ 1.38 million). The business'; C = 'https://www.ejarn.com/article/detail/88973'; D = '2025-06-01' }
    @{ Row = 349; A = 'Join Jamie Kitchen and guests Jörg Saar and John Broughton in the latest episode of "Taking the Temperature on HVACR." Discover '; B = 'Join Jamie Kitchen and guests Jörg Saar and John Broughton in the latest episode of "Taking the Temperature on HVACR." Discover how temperature and pressure interact within refrigeration systems and...'; C = 'https://www.linkedin.com/feed/update/urn:li:activity:7334496373041455105'; D = '2025-05-31' }
    @{ Row = 350; A = 'Danfoss Delivers AC Drives for Finnish Navy’s Four Corvettes'; B = 'Vaasa, Finland-based Danfoss Drives is supplying variable-speed AC drives for the propulsion systems of the Finnish Navy’s four Pohjanmaa-class multi-role corvettes. Part of the Finnish Defense...'; C = 'https://www.maritime-executive.com/corporate/danfoss-delivers-ac-drives-for-finnish-navy-s-four-corvettes'; D = '2025-05-31' }
)

$firstRow = $newRows[0].Row
$lastRow = $newRows[$newRows.Count - 1].Row

# Keep the Date column as plain text (matches the existing A1:D336 data),
# otherwise Excel auto-parses 'yyyy-mm-dd' strings into date serials.
$ws.Range("D" + $firstRow + ":D" + $lastRow).NumberFormat = "@"

foreach ($item in $newRows) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    $ws.Cells.Item($item.Row, 4).Value = $item.D
}

Write-Output ("Wrote rows " + $firstRow + "-" + $lastRow)
